$d = $word.ActiveDocument

# 1. "Multiple agents moving to separate targets..." paragraph: R key now
#    randomises agent positions and targets (previously just "target positions").
$d.Content.Find.ClearFormatting()
$d.Content.Find.Execute("randomises target positions. Removed start", $true, $false, $false, $false, $false, $true, 1, $false, "randomises agent positions and targets. Removed start", 2)

# 2. Insert a new bullet after the "Fixed the bug..." paragraph describing the
#    wall-block rendering / scaling work, keeping the same list formatting.
$d.Content.Find.ClearFormatting()
$d.Content.Find.Execute("agent", $true)
$fixedBugPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "Fixed the bug where sometimes path nodes*") {
        $fixedBugPara = $para
        break
    }
}
$apostrophe = [char]0x2019
$newText = "Changed wall blocks to render as circles, make it look like agents aren" + $apostrophe + "t intersecting with them when moving diagonally. Scaling radius of weapons" + $apostrophe + " effective ranges and agents" + $apostrophe + " radii and avoidance radii according to screen size."
$fixedBugPara.Range.InsertParagraphAfter()
$newParaIndex = $fixedBugPara.Index + 1
$newPara = $d.Paragraphs.Item($newParaIndex)
$newPara.Range.Text = $newText

# 3. "O: toggles highlighting..." paragraph now refers to each agent's optimal
#    path rather than a single optimal path.
$d.Content.Find.ClearFormatting()
$d.Content.Find.Execute(": toggles highlighting of the optimal path in red.", $true, $false, $false, $false, $false, $true, 1, $false, (": toggles highlighting of agents" + $apostrophe + " optimal paths in red."), 2)
